$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.361.31'
$ws.Range("E2").Value = '  -1.25%  '

# Row 3
$ws.Range("D3").Value = '2.188.98'
$ws.Range("E3").Value = '  -1.63%  '

# Row 4
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = "'250.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.14%  '

# Row 6
$ws.Range("D6").Value = "'0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.71%  '

# Row 7
$ws.Range("D7").Value = "'67.13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.05%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = "'0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.87%  '

# Row 10
$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").Value = "'38.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.76%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = "'59.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.23%  '

# Row 12
$ws.Range("D12").Value = "'0.0940"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.99%  '

# Row 13
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("E14").Value = '  -6.02%  '

# Row 15
$ws.Range("D15").Value = '2.516.65'
$ws.Range("E15").Value = '  -1.59%  '

# Row 16
$ws.Range("D16").Value = "'14.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.23%  '

# Row 17
$ws.Range("D17").Value = "'0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.19%  '

# Row 18
$ws.Range("D18").Value = '2.193.88'
$ws.Range("E18").Value = '  -0.82%  '

# Row 19
$ws.Range("D19").Value = '41.307.82'
$ws.Range("E19").Value = '  -1.38%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0948'
$ws.Range("E20").Value = '  -1.90%  '

# Row 21
$ws.Range("D21").Value = "'71.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.36%  '

# Row 22
$ws.Range("D22").Value = "'6.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.12%  '

# Row 23
$ws.Range("D23").Value = "'229.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.51%  '

# Row 24
$ws.Range("D24").Value = "'2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.80%  '

# Row 25
$ws.Range("E25").Value = '  -6.23%  '

# Row 26
$ws.Range("E26").Value = '  +0.13%  '

# Row 27
$ws.Range("D27").Value = "'11.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.53%  '

# Row 28
$ws.Range("D28").Value = "'2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.87%  '

# Row 29
$ws.Range("D29").Value = "'3.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.02%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.16%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = "'166.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.05%  '

# Row 32
$ws.Range("D32").Value = "'20.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.08%  '

# Row 33
$ws.Range("D33").Value = "'0.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.31%  '

# Row 34
$ws.Range("D34").Value = "'5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.90%  '

# Row 35
$ws.Range("D35").Value = "'0.0757"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.04%  '

# Row 36
$ws.Range("D36").Value = "'0.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.26%  '

# Row 37
$ws.Range("D37").Value = "'4.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.79%  '

# Row 38
$ws.Range("D38").Value = "'4.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.19%  '

# Row 39
$ws.Range("D39").Value = "'25.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.55%  '

# Row 40
$ws.Range("D40").Value = "'0.0305"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.82%  '

# Row 41
$ws.Range("E41").Value = '  -3.00%  '

# Row 42
$ws.Range("D42").Value = "'5.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.17%  '

# Row 43
$ws.Range("D43").Value = "'5.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.65%  '

# Row 44
$ws.Range("D44").Value = "'11.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.59%  '

# Row 45
$ws.Range("D45").Value = "'60.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.87%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = "'0.192"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.84%  '

# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'8.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.24%  '

# Row 48
$ws.Range("D48").Value = "'0.0992"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.40%  '

# Row 49
$ws.Range("E49").Value = '  -0.28%  '

# Row 50
$ws.Range("E50").Value = '  -3.18%  '

# Row 51
$ws.Range("D51").Value = "'4.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.78%  '

